# Auto-generated edit script: updates market-data cell values
# (currentAveragePrice / Leve price / profit columns) to match the
# scheduled-runner refresh captured in the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1012.5
$ws.Range("J43").Value = 983.3333
$ws.Range("L43").Value = 983.3333
$ws.Range("N43").Value = -1121.3333
$ws.Range("H62").Value = 33004.39
$ws.Range("I62").Value = 2498.4
$ws.Range("K62").Value = 2498.4
$ws.Range("M62").Value = -1874.4
$ws.Range("H65").Value = 33004.39
$ws.Range("I65").Value = 2498.4
$ws.Range("K65").Value = 12492
$ws.Range("M65").Value = -9372
$ws.Range("H69").Value = 8624.857
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 8624.857
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 25874.571
$ws.Range("N69").Value = -27622.571
$ws.Range("H72").Value = 8624.857
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 8624.857
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 77623.713
$ws.Range("N72").Value = -86359.713
$ws.Range("H98").Value = 2794.2104
$ws.Range("I98").Value = 3490.7693
$ws.Range("J98").Value = 1285
$ws.Range("K98").Value = 3490.7693
$ws.Range("L98").Value = 1285
$ws.Range("M98").Value = -1992.7693
$ws.Range("N98").Value = -4281
$ws.Range("H100").Value = 2618.8333
$ws.Range("I100").Value = 1428.25
$ws.Range("K100").Value = 1428.25
$ws.Range("M100").Value = -887.25
$ws.Range("H112").Value = 1944.4062
$ws.Range("J112").Value = 1860.037
$ws.Range("L112").Value = 5580.111
$ws.Range("N112").Value = -7796.111
$ws.Range("H116").Value = 34396896
$ws.Range("I116").Value = 41861100
$ws.Range("K116").Value = 41861100
$ws.Range("M116").Value = -41857658
$ws.Range("H122").Value = 2794.2104
$ws.Range("I122").Value = 3490.7693
$ws.Range("J122").Value = 1285
$ws.Range("K122").Value = 10472.3079
$ws.Range("L122").Value = 3855
$ws.Range("M122").Value = -8022.3079
$ws.Range("N122").Value = -8755
$ws.Range("H137").Value = 2973
$ws.Range("I137").Value = 2961.8572
$ws.Range("K137").Value = 8885.571599999999
$ws.Range("M137").Value = -6335.571599999999
$ws.Range("H138").Value = 1596.9333
$ws.Range("I138").Value = 1240.85
$ws.Range("K138").Value = 3722.55
$ws.Range("M138").Value = 1417.45
$ws.Range("H141").Value = 5683.4
$ws.Range("I141").Value = 1466
$ws.Range("K141").Value = 4398
$ws.Range("M141").Value = 782

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2329.4285
$ws.Range("I45").Value = 2198.75
$ws.Range("J45").Value = 2747.6
$ws.Range("K45").Value = 2198.75
$ws.Range("L45").Value = 2747.6
$ws.Range("M45").Value = -1821.75
$ws.Range("N45").Value = -3501.6
$ws.Range("H102").Value = 9533.429
$ws.Range("I102").Value = 2952.2666
$ws.Range("J102").Value = 25986.334
$ws.Range("K102").Value = 2952.2666
$ws.Range("L102").Value = 25986.334
$ws.Range("M102").Value = -1330.2666
$ws.Range("N102").Value = -29230.334
$ws.Range("H122").Value = 14495959
$ws.Range("I122").Value = 18521226
$ws.Range("K122").Value = 55563678
$ws.Range("M122").Value = -55561228
$ws.Range("H132").Value = 32260066
$ws.Range("I132").Value = 33335302
$ws.Range("K132").Value = 100005906
$ws.Range("M132").Value = -100003376

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 31939.9
$ws.Range("I20").Value = 42927.145
$ws.Range("J20").Value = 6303
$ws.Range("K20").Value = 42927.145
$ws.Range("L20").Value = 6303
$ws.Range("M20").Value = -42680.145
$ws.Range("N20").Value = -6797
$ws.Range("H94").Value = 2769.3704
$ws.Range("I94").Value = 2762.1904
$ws.Range("J94").Value = 2794.5
$ws.Range("K94").Value = 2762.1904
$ws.Range("L94").Value = 2794.5
$ws.Range("M94").Value = -2311.1904
$ws.Range("N94").Value = -3696.5
$ws.Range("H105").Value = 1625.8572
$ws.Range("I105").Value = 1354.3158
$ws.Range("J105").Value = 2199.111
$ws.Range("K105").Value = 1354.3158
$ws.Range("L105").Value = 2199.111
$ws.Range("M105").Value = 392.6841999999999
$ws.Range("N105").Value = -5693.111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3070.7432
$ws.Range("I31").Value = 1994.2106
$ws.Range("J31").Value = 3442.6365
$ws.Range("K31").Value = 1994.2106
$ws.Range("L31").Value = 3442.6365
$ws.Range("M31").Value = -1699.2106
$ws.Range("N31").Value = -4032.6365
$ws.Range("H34").Value = 3070.7432
$ws.Range("I34").Value = 1994.2106
$ws.Range("J34").Value = 3442.6365
$ws.Range("K34").Value = 1994.2106
$ws.Range("L34").Value = 3442.6365
$ws.Range("M34").Value = -1792.2106
$ws.Range("N34").Value = -3846.6365

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 108.125
$ws.Range("I12").Value = 167
$ws.Range("J12").Value = 94.53846
$ws.Range("K12").Value = 501
$ws.Range("L12").Value = 283.61538
$ws.Range("M12").Value = -328
$ws.Range("N12").Value = -629.61538
$ws.Range("H22").Value = 316.33334
$ws.Range("I22").Value = 316.33334
$ws.Range("K22").Value = 949.0000200000001
$ws.Range("M22").Value = -780.0000200000001
$ws.Range("H27").Value = 316.33334
$ws.Range("I27").Value = 316.33334
$ws.Range("K27").Value = 949.0000200000001
$ws.Range("M27").Value = -847.0000200000001
$ws.Range("H34").Value = 815.1539
$ws.Range("I34").Value = 606
$ws.Range("J34").Value = 945.875
$ws.Range("K34").Value = 1818
$ws.Range("L34").Value = 2837.625
$ws.Range("M34").Value = -1734
$ws.Range("N34").Value = -3005.625
$ws.Range("H38").Value = 212.63637
$ws.Range("I38").Value = 280.25
$ws.Range("J38").Value = 32.333332
$ws.Range("K38").Value = 840.75
$ws.Range("L38").Value = 96.999996
$ws.Range("M38").Value = -493.75
$ws.Range("N38").Value = -790.999996
$ws.Range("H50").Value = 350
$ws.Range("J50").Value = 350
$ws.Range("L50").Value = 1050
$ws.Range("N50").Value = -2012
$ws.Range("H53").Value = 350
$ws.Range("J53").Value = 350
$ws.Range("L53").Value = 1050
$ws.Range("N53").Value = -2012
$ws.Range("H68").Value = 1683
$ws.Range("I68").Value = 1574.75
$ws.Range("J68").Value = 1899.5
$ws.Range("K68").Value = 4724.25
$ws.Range("L68").Value = 5698.5
$ws.Range("M68").Value = -3913.25
$ws.Range("N68").Value = -7320.5
$ws.Range("H71").Value = 1683
$ws.Range("I71").Value = 1574.75
$ws.Range("J71").Value = 1899.5
$ws.Range("K71").Value = 14172.75
$ws.Range("L71").Value = 17095.5
$ws.Range("M71").Value = -10116.75
$ws.Range("N71").Value = -25207.5
$ws.Range("H113").Value = 1689.1111
$ws.Range("I113").Value = 675
$ws.Range("K113").Value = 2025
$ws.Range("M113").Value = 145
$ws.Range("H122").Value = 596.1667
$ws.Range("I122").Value = 566.3333
$ws.Range("J122").Value = 626
$ws.Range("K122").Value = 5096.9997
$ws.Range("L122").Value = 5634
$ws.Range("M122").Value = -2646.9997
$ws.Range("N122").Value = -10534
$ws.Range("H131").Value = 8496.615
$ws.Range("I131").Value = 971.5
$ws.Range("K131").Value = 2914.5
$ws.Range("M131").Value = 2125.5
$ws.Range("H132").Value = 2228.3333
$ws.Range("J132").Value = 2228.3333
$ws.Range("L132").Value = 20054.9997
$ws.Range("N132").Value = -25114.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 11999.208
$ws.Range("J70").Value = 12332.75
$ws.Range("L70").Value = 12332.75
$ws.Range("N70").Value = -12872.75
$ws.Range("H73").Value = 11999.208
$ws.Range("J73").Value = 12332.75
$ws.Range("L73").Value = 12332.75
$ws.Range("N73").Value = -14204.75
$ws.Range("H102").Value = 2053.5667
$ws.Range("I102").Value = 1130.2222
$ws.Range("K102").Value = 1130.2222
$ws.Range("M102").Value = 491.7778000000001
$ws.Range("H126").Value = 8189.75
$ws.Range("J126").Value = 4545.8887
$ws.Range("L126").Value = 13637.6661
$ws.Range("N126").Value = -18577.6661

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I7").Value = 17243738
$ws.Range("J7").Value = 3639.9092
$ws.Range("K7").Value = 17243738
$ws.Range("L7").Value = 3639.9092
$ws.Range("M7").Value = -17243626
$ws.Range("N7").Value = -3863.9092
$ws.Range("H22").Value = 659.75
$ws.Range("J22").Value = 616.7143
$ws.Range("L22").Value = 616.7143
$ws.Range("N22").Value = -1206.7143
$ws.Range("H27").Value = 659.75
$ws.Range("J27").Value = 616.7143
$ws.Range("L27").Value = 616.7143
$ws.Range("N27").Value = -830.7143
$ws.Range("H43").Value = 5026311.5
$ws.Range("I43").Value = 3595571.5
$ws.Range("K43").Value = 3595571.5
$ws.Range("M43").Value = -3595378.5
$ws.Range("H55").Value = 516.05884
$ws.Range("I55").Value = 388.8889
$ws.Range("K55").Value = 388.8889
$ws.Range("M55").Value = -215.8889
$ws.Range("I126").Value = 17243738
$ws.Range("J126").Value = 3639.9092
$ws.Range("K126").Value = 51731214
$ws.Range("L126").Value = 10919.7276
$ws.Range("M126").Value = -51728744

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("M69").ClearContents()
$ws.Range("M72").ClearContents()
